# Change russian text to svg_0
#
# Six rows on the "Translation" sheet referenced Text widgets whose text was
# Russian (Логин, Пароль, А.Вход_1: .. А.Вход_4:). Those widgets were
# converted to SVG images (svg_0, ...), so they are no longer text widgets
# and must be removed from the generated Translation table. Deleting the
# rows shifts everything below them up, which is exactly what the diff
# shows (uniqueCount of the shared strings table is unchanged, only row
# contents shift and the six rows that used to be at the bottom become
# blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Original row numbers (before any deletion) holding the Russian-language
# rows that need to disappear.
$rowsToDelete = @(11, 15, 19, 20, 21, 22)

# Delete from the bottom up so earlier row numbers in the list stay valid.
$sortedRows = $rowsToDelete | Sort-Object -Descending
foreach ($r in $sortedRows) {
    $ws.Rows("$r`:$r").Delete()
}
